# Add new column 'Correction ' to Card1 by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# Header cell N1: same text + same formatting as the rest of the header row
# (bold / bordered / centred), copied from M1 (the previous last column).
$ws.Cells.Item(1, 14).Value = "Correction "
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the previously-empty M column (rows 2-12) with "nan" to match the rest
# of the column, and create the new (empty) N column cells for the same
# rows so the sheet's used range grows to N12, same as the source edit.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"   # column M

    $ws.Cells.Item($r, 14).Value = ""      # column N (blank)
    $ws.Cells.Item($r, 14).Style = "Normal"
}
